$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 662.13794
$ws.Range("I33").Value = 342.72
$ws.Range("K33").Value = 342.72
$ws.Range("M33").Value = -113.72
$ws.Range("H98").Value = 920.90625
$ws.Range("I98").Value = 898.9666999999999
$ws.Range("J98").Value = 1250
$ws.Range("K98").Value = 898.9666999999999
$ws.Range("L98").Value = 1250
$ws.Range("M98").Value = 599.0333000000001
$ws.Range("N98").Value = -4246
$ws.Range("H112").Value = 1828.1818
$ws.Range("I112").Value = 866.6667
$ws.Range("J112").Value = 1980
$ws.Range("K112").Value = 2600.0001
$ws.Range("L112").Value = 5940
$ws.Range("M112").Value = -1492.0001
$ws.Range("N112").Value = -8156
$ws.Range("H122").Value = 920.90625
$ws.Range("I122").Value = 898.9666999999999
$ws.Range("J122").Value = 1250
$ws.Range("K122").Value = 2696.9001
$ws.Range("L122").Value = 3750
$ws.Range("M122").Value = -246.9000999999998
$ws.Range("N122").Value = -8650
$ws.Range("H132").Value = 6101016
$ws.Range("I132").Value = 3495.7576
$ws.Range("J132").Value = 31253288
$ws.Range("K132").Value = 10487.2728
$ws.Range("L132").Value = 93759864
$ws.Range("M132").Value = -7957.272799999999
$ws.Range("N132").Value = -93764924
$ws.Range("H141").Value = 1081.8485
$ws.Range("I141").Value = 990.35486
$ws.Range("J141").Value = 2500
$ws.Range("K141").Value = 2971.06458
$ws.Range("L141").Value = 7500
$ws.Range("M141").Value = 2208.93542
$ws.Range("N141").Value = -17860

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6946.09
$ws.Range("I32").Value = 5529.427
$ws.Range("J32").Value = 18408.182
$ws.Range("K32").Value = 5529.427
$ws.Range("L32").Value = 18408.182
$ws.Range("M32").Value = -5242.427
$ws.Range("N32").Value = -18982.182
$ws.Range("H45").Value = 2166142.5
$ws.Range("I45").Value = 2842524.5
$ws.Range("J45").Value = 1720
$ws.Range("K45").Value = 2842524.5
$ws.Range("L45").Value = 1720
$ws.Range("M45").Value = -2842147.5
$ws.Range("N45").Value = -2474
$ws.Range("H61").Value = 2206.1633
$ws.Range("I61").Value = 2229.7273
$ws.Range("J61").Value = 1998.8
$ws.Range("K61").Value = 2229.7273
$ws.Range("L61").Value = 1998.8
$ws.Range("M61").Value = -2017.7273
$ws.Range("N61").Value = -2422.8
$ws.Range("H110").Value = 1179.85
$ws.Range("I110").Value = 874.8125
$ws.Range("J110").Value = 2400
$ws.Range("K110").Value = 874.8125
$ws.Range("L110").Value = 2400
$ws.Range("M110").Value = 1170.1875
$ws.Range("N110").Value = -6490
$ws.Range("H136").Value = 2206.1633
$ws.Range("I136").Value = 2229.7273
$ws.Range("J136").Value = 1998.8
$ws.Range("K136").Value = 6689.1819
$ws.Range("L136").Value = 5996.4
$ws.Range("M136").Value = -4139.1819
$ws.Range("N136").Value = -11096.4

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4347.5864
$ws.Range("I105").Value = 2900
$ws.Range("J105").Value = 4454.815
$ws.Range("K105").Value = 2900
$ws.Range("L105").Value = 4454.815
$ws.Range("M105").Value = -1153
$ws.Range("N105").Value = -7948.815

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3190.362
$ws.Range("I31").Value = 2992.818
$ws.Range("J31").Value = 3311.0833
$ws.Range("K31").Value = 2992.818
$ws.Range("L31").Value = 3311.0833
$ws.Range("M31").Value = -2697.818
$ws.Range("N31").Value = -3901.0833
$ws.Range("H34").Value = 3190.362
$ws.Range("I34").Value = 2992.818
$ws.Range("J34").Value = 3311.0833
$ws.Range("K34").Value = 2992.818
$ws.Range("L34").Value = 3311.0833
$ws.Range("M34").Value = -2790.818
$ws.Range("N34").Value = -3715.0833
$ws.Range("H99").Value = 1902.8518
$ws.Range("I99").Value = 1611.5
$ws.Range("J99").Value = 2326.6365
$ws.Range("K99").Value = 1611.5
$ws.Range("L99").Value = 2326.6365
$ws.Range("M99").Value = -113.5
$ws.Range("N99").Value = -5322.636500000001
$ws.Range("H126").Value = 1902.8518
$ws.Range("I126").Value = 1611.5
$ws.Range("J126").Value = 2326.6365
$ws.Range("K126").Value = 4834.5
$ws.Range("L126").Value = 6979.9095
$ws.Range("M126").Value = -2364.5
$ws.Range("N126").Value = -11919.9095
$ws.Range("H132").Value = 2705.1738
$ws.Range("I132").Value = 2255.6875
$ws.Range("J132").Value = 3732.5715
$ws.Range("K132").Value = 6767.0625
$ws.Range("L132").Value = 11197.7145
$ws.Range("M132").Value = -4237.0625
$ws.Range("N132").Value = -16257.7145

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 694549.2
$ws.Range("I2").Value = 72.545456
$ws.Range("J2").Value = 1543354
$ws.Range("K2").Value = 435.272736
$ws.Range("L2").Value = 9260124
$ws.Range("M2").Value = -322.272736
$ws.Range("N2").Value = -9260350
$ws.Range("H113").Value = 58824364
$ws.Range("I113").Value = 166667280
$ws.Range("J113").Value = 959.0909
$ws.Range("K113").Value = 500001840
$ws.Range("L113").Value = 2877.2727
$ws.Range("M113").Value = -499999670
$ws.Range("N113").Value = -7217.2727
$ws.Range("H131").Value = 1895.4459
$ws.Range("I131").Value = 3318.6
$ws.Range("J131").Value = 1533.6271
$ws.Range("K131").Value = 9955.799999999999
$ws.Range("L131").Value = 4600.8813
$ws.Range("M131").Value = -4915.799999999999
$ws.Range("N131").Value = -14680.8813
$ws.Range("H134").Value = 4788
$ws.Range("I134").Value = 2520
$ws.Range("J134").Value = 8190
$ws.Range("K134").Value = 7560
$ws.Range("L134").Value = 24570
$ws.Range("M134").Value = -2490
$ws.Range("N134").Value = -34710
$ws.Range("H139").Value = 3343.75
$ws.Range("I139").Value = 1750
$ws.Range("J139").Value = 4937.5
$ws.Range("K139").Value = 5250
$ws.Range("L139").Value = 14812.5
$ws.Range("M139").Value = -110
$ws.Range("N139").Value = -25092.5
$ws.Range("H140").Value = 3745.077
$ws.Range("I140").Value = 1698.7273
$ws.Range("J140").Value = 15000
$ws.Range("K140").Value = 5096.1819
$ws.Range("L140").Value = 45000
$ws.Range("M140").Value = 83.81810000000041
$ws.Range("N140").Value = -55360

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 32000
$ws.Range("I70").Value = 102000
$ws.Range("K70").Value = 102000
$ws.Range("M70").Value = -101730
$ws.Range("H73").Value = 32000
$ws.Range("I73").Value = 102000
$ws.Range("K73").Value = 102000
$ws.Range("M73").Value = -101064

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 249.55
$ws.Range("I55").Value = 233.875
$ws.Range("J55").Value = 260
$ws.Range("K55").Value = 233.875
$ws.Range("L55").Value = 260
$ws.Range("M55").Value = -60.875
$ws.Range("N55").Value = -606

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 42946.11
$ws.Range("J135").Value = 42946.11
$ws.Range("L135").Value = 42946.11
$ws.Range("N135").Value = -53086.11
